# Update Name of Algo
# Apply updated imputed values to result_data_KNN sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = -8.436
$ws.Range("D6").Value = -8.055999999999999
$ws.Range("C7").Value = -13.498
$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.808
$ws.Range("E12").Value = 13.041
$ws.Range("B13").Value = 6.955
$ws.Range("A18").Value = -21.78
$ws.Range("C20").Value = -12.879
$ws.Range("E20").Value = 12.637
$ws.Range("E25").Value = 12.784
